# Apply the "update new orleans xlsx files" edit:
#  - hotel_info sheet gains a new "State" column (inserted after Hotel_Name, before City)
#    and the data row is updated accordingly.
#  - review_info sheet becomes the first tab, hotel_info becomes the second tab.

$wb = $excel.ActiveWorkbook

$hotel = $wb.Worksheets.Item("hotel_info")
$review = $wb.Worksheets.Item("review_info")

# --- hotel_info: insert a new "State" column before "City" (column C) ---
$hotel.Columns.Item(3).Insert()

$hotel.Range("C1").Value = "State"

# --- hotel_info: write out the data row exactly as it ends up after the edit ---
$hotel.Range("A2").Value = 43024
$hotel.Range("B2").Value = "Total_Reviews_num"
$hotel.Range("C2").Value = "French Quarter Guest Houses"
$hotel.Range("D2").Value = "Louisiana"
$hotel.Range("E2").Value = 70116
$hotel.Range("F2").Value = "??? Couldn't find"

# --- reorder tabs: review_info first, hotel_info second ---
$review.Move($hotel)
